$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap cell formatting (font color / style) between B and C
#     for the rows where the diff shows the style indices exchanged.
$swapRows = 3,6,14,22,26
foreach ($r in $swapRows) {
    $bAddr = "B$r"
    $cAddr = "C$r"
    $stageAddr = "Z$r"
    $ws.Range($bAddr).Copy()
    $ws.Range($stageAddr).PasteSpecial(-4122)
    $ws.Range($cAddr).Copy()
    $ws.Range($bAddr).PasteSpecial(-4122)
    $ws.Range($stageAddr).Copy()
    $ws.Range($cAddr).PasteSpecial(-4122)
    $ws.Range($stageAddr).Clear()
}
$excel.CutCopyMode = $false

# --- Step 2: update the cell values in column B (and C2) to the new numbers
$ws.Range("B2").Value = [double]"16.0"
$ws.Range("B3").Value = [double]"110689095.0"
$ws.Range("B4").Value = [double]"1.670124101e+33"
$ws.Range("B5").Value = [double]"224945.43"
$ws.Range("B6").Value = [double]"196.99525"
$ws.Range("B7").Value = [double]"174.98017"
$ws.Range("B8").Value = [double]"2.7551077"
$ws.Range("B9").Value = [double]"330.411"
$ws.Range("B10").Value = [double]"180.96006"
$ws.Range("B11").Value = [double]"2479.0088"
$ws.Range("B12").Value = [double]"4651.2664"
$ws.Range("B13").Value = [double]"7980.7025"
$ws.Range("B14").Value = [double]"54612870.0"
$ws.Range("B15").Value = [double]"32233928.5"
$ws.Range("B16").Value = [double]"1615669.87"
$ws.Range("B17").Value = [double]"13321995.4"
$ws.Range("B18").Value = [double]"1671.8404"
$ws.Range("B19").Value = [double]"1058.2745"
$ws.Range("B20").Value = [double]"4282734.4"
$ws.Range("B21").Value = [double]"3208229.3"
$ws.Range("B22").Value = [double]"917.78449"
$ws.Range("B23").Value = [double]"358.70562"
$ws.Range("B24").Value = [double]"3482.64756"
$ws.Range("B25").Value = [double]"632.7367"
$ws.Range("B26").Value = [double]"853.23778"
$ws.Range("B27").Value = [double]"590.64621"
$ws.Range("B28").Value = [double]"675.4623"
$ws.Range("B29").Value = [double]"675.30388"
$ws.Range("B30").Value = [double]"577.35496"
$ws.Range("B31").Value = [double]"1070.16596"
$ws.Range("B32").Value = [double]"5917658.0"
$ws.Range("C2").Value = [double]"13.0"
